$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Decrement by 1 the D column values for "farms_total_count" (row has C = "farms_total_count")
# and "farms_to_examine_count" rows, for every existing week row 2..68 (rows where
# column C is farms_total_count or farms_to_examine_count), matching original values.
$rowsToDecrement = @(2,3,7,8,12,13,17,18,22,23,27,28,32,33,37,38,42,43,47,48,52,53,57,58,62,63,67,68)
foreach ($r in $rowsToDecrement) {
    $cell = $ws.Cells.Item($r, 4)
    $current = $cell.Value()
    $cell.Value = $current - 1
}

# Append new rows 72-76 for YearWeekIso 202502 (LastDayOfWeek 2025-01-12, serial 45669)
$newRows = @(
    @{Row=72; Variable="farms_total_count"; Number=12361},
    @{Row=73; Variable="farms_to_examine_count"; Number=5130},
    @{Row=74; Variable="farms_examined_count"; Number=7231},
    @{Row=75; Variable="farms_examined_positive_count"; Number=1417},
    @{Row=76; Variable="farms_examined_negative_count"; Number=5814}
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value = 202502
    $dateCell = $ws.Cells.Item($r, 2)
    # Reuse the same date style already used by the other LastDayOfWeek cells
    # (numFmtId 14, "m/d/yyyy") instead of creating a brand-new custom
    # number format by assigning .NumberFormat directly.
    $ws.Cells.Item(67, 2).Copy()
    $dateCell.PasteSpecial(-4122)
    $dateCell.Value = 45669
    $ws.Cells.Item($r, 3).Value = $nr.Variable
    $ws.Cells.Item($r, 4).Value = $nr.Number
}

$ws.Range("G14").Select()
